$wb = $excel.ActiveWorkbook

# --- Update "jobs" sheet due dates (C2:C11) ---
$jobs = $wb.Worksheets.Item("jobs")
$jobs.Range("C2").Value = 43746.75
$jobs.Range("C3").Value = 43747.75
$jobs.Range("C4").Value = 43745.75
$jobs.Range("C5").Value = 43746.75
$jobs.Range("C6").Value = 43745.75
$jobs.Range("C7").Value = 43746.75
$jobs.Range("C8").Value = 43746.75
$jobs.Range("C9").Value = 43745.75
$jobs.Range("C10").Value = 43745.75
$jobs.Range("C11").Value = 43745.75

# --- Update selection/active sheet ---
# "machines" sheet was previously the active/selected tab; now "jobs" is active
# with selection at C18.
$jobs.Select()
$jobs.Range("C18").Select()

$machines = $wb.Worksheets.Item("machines")
$machines.Range("H8").Select()

$jobs.Select()
